# New PO forecast model
# Updates the three PO-analysis sheets with the latest weekly/monthly POs
# and a re-fit forecast curve.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Weekly Quantity": append the two most recent weekly observations.
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

$wsWeekly.Range("A55").Value = 45676.99999999999
$wsWeekly.Range("A55").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWeekly.Range("B55").Value = 5

$wsWeekly.Range("A56").Value = 45683.99999999999
$wsWeekly.Range("A56").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWeekly.Range("B56").Value = 1

# ---------------------------------------------------------------------
# "Monthly Trend": append the newest monthly roll-up.
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsMonthly.Range("A21").Value = 45688.99999999999
$wsMonthly.Range("A21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsMonthly.Range("B21").Value = 6

# ---------------------------------------------------------------------
# "PO Forecast": refit forecast values for existing weeks and extend the
# forecast horizon with new weeks.
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

$wsForecast.Range("B2:B4").Value = 7
$wsForecast.Range("B10:B39").Value = 8
$wsForecast.Range("B40:B54").Value = 9

# Rows 55-62 move onto the new weekly cadence (new dates) at the new
# forecast level.
$newForecastDates = @(
    45676.99999999999,
    45683.99999999999,
    45690.99999999999,
    45697.99999999999,
    45704.99999999999,
    45711.99999999999,
    45718.99999999999,
    45725.99999999999,
    45732.99999999999,
    45739.99999999999
)

$row = 55
foreach ($d in $newForecastDates) {
    $cellA = $wsForecast.Cells.Item($row, 1)
    $cellA.Value = $d
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($row, 2).Value = 10
    $row++
}
